# Apply cryptos list update (values scraped Sat Jul  1 19:34:20 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Text" number format on Price cells whose new values look numeric,
# so Excel keeps them as literal strings (matching the source data) instead
# of auto-converting them to numbers.
$textCells = @(
    "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D23", "D24", "D25", "D26", "D27", "D31", "D32", "D33", "D34", "D36", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "30.627.77"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.923.35"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "247.24"
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.4745"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.2889"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").Value = "0.06841"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").Value = "105.18"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").Value = "18.38"
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("D12").Value = "1.921.73"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "0.07698"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "5.323"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").Value = "0.6681"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "291.63"
$ws.Range("E16").Value = "  -2.73%  "
$ws.Range("D17").Value = "30.636.19"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000007627"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "5.590"
$ws.Range("E19").Value = "  +7.07%  "
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "2.167.88"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "6.437"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").Value = "9.454"
$ws.Range("E25").Value = "  +3.07%  "
$ws.Range("D26").Value = "167.86"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "21.10"
$ws.Range("E28").Value = "  +5.19%  "
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").Value = "4.184"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "4.065"
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("D33").Value = "0.05044"
$ws.Range("D34").Value = "0.7387"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "0.02073"
$ws.Range("E36").Value = "  +6.53%  "
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "2.062"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("D41").Value = "0.8775"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Value = "0.4391"
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("D43").Value = "5.902"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "68.08"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").Value = "7.277"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "9.313"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").Value = "48.24"
$ws.Range("E48").Value = "  +15.02%  "
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("D50").Value = "0.2523"
$ws.Range("E50").Value = "  +13.55%  "
$ws.Range("D51").Value = "35.04"
$ws.Range("E51").Value = "  +1.33%  "
